$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 6 (ano = 2025) values
$ws.Range("C6").Value = 433
$ws.Range("E6").Value = 124
$ws.Range("G6").Value = 28.63741339491917
$ws.Range("H6").Value = 71.36258660508084
